$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, reusing the same style as the other
# header cells (e.g. G1) by copying the formatted cell and then overwriting
# its value/text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add new value in H2 (plain numeric cell, no special style)
$ws.Range("H2").Value = 1
